# Revert "Revert "added ZX correction""
#
# 1) RSD_Models sheet: row 4 changes from "Zephir" to "ZX" (the model list
#    entry is swapped out for the ZX correction device).
# 2) Configuration_template sheet: the RSD_model selector (B34) and the
#    correction selector (B36) are updated to reflect the new pick
#    ("Triton" / "ZX").
# 3) View-state touch-ups that came along with the edit (selection on the
#    RSD_Models sheet, scroll position on Configuration_template, and an
#    autofit of RSD_Models column A).

$wb = $excel.ActiveWorkbook

$wsConfig = $wb.Worksheets.Item("Configuration_template")
$wsModels = $wb.Worksheets.Item("RSD_Models")

# --- RSD_Models: "Zephir" -> "ZX" -------------------------------------
$wsModels.Range("A4").Value = "ZX"

# --- Configuration_template: update the two dropdown-style picks -------
$wsConfig.Range("B34").Value = "Triton"
$wsConfig.Range("B36").Value = "ZX"

# --- View state: RSD_Models selection moves to A5, column A autofits ---
$wsModels.Activate()
$wsModels.Columns("A:A").AutoFit()
$wsModels.Range("A5").Select()

# --- View state: Configuration_template scrolls so row 19 is on top ----
$wsConfig.Activate()
$winConfig = $excel.ActiveWindow
$winConfig.ScrollRow = 19
$winConfig.ScrollColumn = 1
